$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Status"

$ws.Range("B2").Value = "Child"
$ws.Range("C2").Value = "Free"
$ws.Range("D2").Value = "Y"

$ws.Range("B3").Value = "Adult"
$ws.Range("C3").Value = "Free"
$ws.Range("D3").Value = "Y"

$ws.Range("B4").Value = "Child"
$ws.Range("C4").Value = "Free"
$ws.Range("D4").Value = "Y"

$ws.Range("B5").Value = "Student"
$ws.Range("C5").Value = "Free"
$ws.Range("D5").Value = "Y"

$ws.Range("B6").Value = "Adult"
$ws.Range("C6").Value = "Linked"

$ws.Range("B7").Value = "Adult"
$ws.Range("C7").Value = "Free"
$ws.Range("D7").Value = "Y"

$ws.Range("B8").Value = "Adult"
$ws.Range("C8").Value = "Linked"

$ws.Range("B9").Value = "Student"
$ws.Range("C9").Value = "Free"
$ws.Range("D9").Value = "Y"

$ws.Range("B10").Value = "Adult"
$ws.Range("C10").Value = "Free"
$ws.Range("D10").Value = "Y"

$ws.Range("B11").Value = "Child"
$ws.Range("C11").Value = "Free"
$ws.Range("D11").Value = "Y"

$ws.Range("D11").Select()
